$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("interface_calibration")

# Update the AVERAGE formula range in E12 (now includes E3:E10 instead of E6:E10)
$ws.Range("E12").Formula = "=AVERAGE(E3:E10)"

# Update the measured/weighed value in C15 and apply a thousands-separator number format
$ws.Range("C15").Value = 265.55983950000001
$ws.Range("C15").NumberFormat = "#,##0"

# Update the selected cell to reflect where the user last clicked
$ws.Range("C16").Select()

$wb.Application.Calculate()
